$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Nombre completo" header in column C --------------------------------
$ws.Range("C1").Value = "Nombre completo"

# --- Insert the new "Skyway" customer row right after the header -------------
$ws.Rows.Item(2).Insert()
$ws.Range("A2:C2").ClearFormats()
$ws.Range("A2").Value = "Skyway"
$ws.Range("B2").Value = 7001353
$ws.Range("C2").Value = "Skyway"

# --- Convert the sap_code values that were stored as text into real numbers --
# (existing rows shifted down by one after the insert above)
$ws.Range("B4").Value = 70004530   # Coffeyville Sektam
$ws.Range("B5").Value = 70016983   # Concentric India
$ws.Range("B6").Value = 70001353   # Daleo
$ws.Range("B11").Value = 70017109  # Hypro
$ws.Range("B13").Value = 70017128  # PARKER
$ws.Range("B16").Value = 70013219  # Sauer-Danfoss
$ws.Range("B17").Value = 70018728  # Soucy

# --- Sort the whole table (header included) by the sap_code column -----------
$sortRange = $ws.Range("A1:C25")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B1"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Turn on the AutoFilter for the header row --------------------------------
$ws.Range("A1:C1").AutoFilter()

# --- Register the hidden _FilterDatabase defined name (mirrors Excel's own
#     bookkeeping whenever AutoFilter is toggled on from the ribbon) ----------
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$1")
$fd.Visible = $false

# --- Restore the selected cell as recorded in the saved file ------------------
$ws.Range("D9").Select()
